$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '82.120.04'
$ws.Range('E2').Value = '  +3.13%  '
$ws.Range('D3').Value = '3.165.94'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '219.00'
$ws.Range('E5').Value = '  +5.85%  '
$ws.Range('D6').Value = '619.24'
$ws.Range('E6').Value = '  -3.27%  '
$ws.Range('D7').Value = '0.291'
$ws.Range('E7').Value = '  +18.96%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.582'
$ws.Range('E9').Value = '  -4.34%  '
$ws.Range('D10').Value = '3.164.39'
$ws.Range('E10').Value = '  -1.08%  '
$ws.Range('D11').Value = '0.600'
$ws.Range('E11').Value = '  -2.92%  '
$ws.Range('D12').Value = '0.0000257'
$ws.Range('E12').Value = '  -1.00%  '
$ws.Range('E13').Value = '  -0.68%  '
$ws.Range('D14').Value = '5.31'
$ws.Range('E14').Value = '  -2.64%  '
$ws.Range('D15').Value = '3.740.87'
$ws.Range('E15').Value = '  -1.38%  '
$ws.Range('D16').Value = '32.28'
$ws.Range('E16').Value = '  -1.16%  '
$ws.Range('D17').Value = '81.855.23'
$ws.Range('E17').Value = '  +2.91%  '
$ws.Range('D18').Value = '3.156.26'
$ws.Range('E18').Value = '  -1.49%  '
$ws.Range('D19').Value = '3.25'
$ws.Range('E19').Value = '  +9.28%  '
$ws.Range('D20').Value = '14.00'
$ws.Range('E20').Value = '  -4.74%  '
$ws.Range('D21').Value = '435.57'
$ws.Range('E21').Value = '  -2.01%  '
$ws.Range('D22').Value = '8.91'
$ws.Range('E22').Value = '  -5.95%  '
$ws.Range('D23').Value = '5.14'
$ws.Range('E23').Value = '  -3.75%  '
$ws.Range('D24').Value = '7.26'
$ws.Range('E24').Value = '  +4.06%  '
$ws.Range('D25').Value = '5.26'
$ws.Range('E25').Value = '  +8.78%  '
$ws.Range('D26').Value = '11.90'
$ws.Range('E26').Value = '  +8.60%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').Value = '76.74'
$ws.Range('E27').Value = '  -1.43%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '3.313.60'
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  -1.12%  '
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').Value = '9.03'
$ws.Range('E32').Value = '  -2.21%  '
$ws.Range('D33').Value = '570.67'
$ws.Range('E33').Value = '  +5.44%  '
$ws.Range('E34').Value = '  -2.32%  '
$ws.Range('D35').Value = '0.148'
$ws.Range('E35').Value = '  +20.01%  '
$ws.Range('D36').Value = '0.152'
$ws.Range('E36').Value = '  -0.56%  '
$ws.Range('D37').Value = '1.99'
$ws.Range('E37').Value = '  -2.54%  '
$ws.Range('D38').Value = '22.69'
$ws.Range('E38').Value = '  -2.91%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').Value = '6.23'
$ws.Range('E39').Value = '  +10.99%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').Value = '0.406'
$ws.Range('E41').Value = '  -1.85%  '
$ws.Range('D42').Value = '20.85'
$ws.Range('E42').Value = '  +4.04%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '2.02'
$ws.Range('E43').Value = '  +10.19%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '3.02'
$ws.Range('E44').Value = '  +13.60%  '
$ws.Range('D45').Value = '158.86'
$ws.Range('E45').Value = '  -3.93%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').Value = '186.58'
$ws.Range('E47').Value = '  -4.23%  '
$ws.Range('D48').Value = '44.56'
$ws.Range('E48').Value = '  +1.50%  '
$ws.Range('D50').Value = '26.33'
$ws.Range('E50').Value = '  +1.45%  '
$ws.Range('D51').Value = '0.766'
$ws.Range('E51').Value = '  -5.42%  '
